# "Completed prototype of foundations chapter"
#
# Content changes (everything else in the two touched paragraphs/runs is a
# byte-for-byte re-split of already-existing text, i.e. no wording changes):
#
#  - Slide 12 ("Quellen (1)"): the "Encrypted " run in the [8] citation title
#    is split into two runs, "Encrypted" + " ".
#  - Slide 5 ("CryptDB"):
#      * "Proxy zwischen Client und Server übersetzt" is split into
#        "Proxy zwischen Client und Server " + "übersetzt".
#      * "Anfragen des Client in verschlüsselte Form" becomes
#        "Anfragen des Client in optimierte Form" (verschlüsselte -> optimierte),
#        split into "Anfragen des Client " + "in optimierte Form".
#      * "Spezielles Verschlüsselungsschema (SQL-aware " is split into
#        "Spezielles " + "Verschlüsselungsschema (SQL-aware ".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 12 - "Quellen (1)" - content placeholder - paragraph 8 ([8] citation)
# ---------------------------------------------------------------------------
$slide12 = $p.Slides.Item(12)
$shape12 = $slide12.Shapes.Item("Inhaltsplatzhalter 2")
$tr12 = $shape12.TextFrame.TextRange
$para8 = $tr12.Paragraphs(8, 1)

$text8 = $para8.Text
$encIdx = $text8.IndexOf("Encrypted ")
if ($encIdx -lt 0) { throw "Could not locate 'Encrypted ' run on slide 12" }
# "Encrypted" is 9 characters long; re-writing just that prefix splits the
# run right after it, so the trailing space becomes its own run.
$encRun = $para8.Characters($encIdx + 1, 9)
$encRun.Text = "Encrypted"

# ---------------------------------------------------------------------------
# Slide 5 - "CryptDB" - content placeholder
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$shape5 = $slide5.Shapes.Item("Inhaltsplatzhalter 2")
$tr5 = $shape5.TextFrame.TextRange

# Paragraph 3 (lvl 1): "Proxy zwischen Client und Server übersetzt"
#   -> "Proxy zwischen Client und Server " + "übersetzt"
$para3 = $tr5.Paragraphs(3, 1)
$text3 = $para3.Text
$wordIdx3 = $text3.IndexOf("übersetzt")
if ($wordIdx3 -lt 0) { throw "Could not locate 'uebersetzt' in slide 5 paragraph 3" }
$splitRun3 = $para3.Characters($wordIdx3 + 1, $text3.Length - $wordIdx3)
$splitRun3.Text = "übersetzt"

# Paragraph 4 (lvl 2): "Anfragen des Client in verschlüsselte Form"
#   -> "Anfragen des Client " + "in optimierte Form"
$para4 = $tr5.Paragraphs(4, 1)
$text4 = $para4.Text
$wordIdx4 = $text4.IndexOf("in verschl")
if ($wordIdx4 -lt 0) { throw "Could not locate 'in verschluesselte' in slide 5 paragraph 4" }
$splitRun4 = $para4.Characters($wordIdx4 + 1, $text4.Length - $wordIdx4)
$splitRun4.Text = "in optimierte Form"

# Paragraph 6 (lvl 1): "Spezielles Verschlüsselungsschema (SQL-aware " + ...
#   -> "Spezielles " + "Verschlüsselungsschema (SQL-aware " + ...
# (only the leading run is touched; "adaptable", " ", "encryption", etc.
#  that follow it in the same paragraph are left untouched)
$para6 = $tr5.Paragraphs(6, 1)
$text6 = $para6.Text
$wordIdx6 = $text6.IndexOf("Verschl")
$runEndIdx6 = $text6.IndexOf("adaptable")
if ($wordIdx6 -lt 0 -or $runEndIdx6 -lt 0) { throw "Could not locate split points in slide 5 paragraph 6" }
$splitRun6 = $para6.Characters($wordIdx6 + 1, $runEndIdx6 - $wordIdx6)
$splitRun6.Text = $splitRun6.Text
